$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- C2: updated TPL subject/note text ---
$c2Text = @"
Subject: Theory of Programming Languages;
Instructor: Dr Amjad Mehmood;
ClassSenior: +92 301 3068788‬;
Note: First 5 Chapters are included in Mid Term;
CreditHours: 3.0;
"@
$ws.Range("C2").Value = $c2Text

# --- C3: Week 1 topic (new content) ---
$c3Text = @"
Topic: name- Week 1, lectures- 3 Lectures, duration- 01:47;
Video: link- https://drive.google.com/file/d/1ctSBf7a08q8Au16-1HEayw8aghDsauQG/preview, name- TPL Week#1 Part1, duration- 00:17;
Video: link- https://drive.google.com/file/d/17yY-6RGVlhAsfMqJVqewvXUHGmNDlcR1/preview, name- TPL Week#1 Part2, duration- 00:53;
Slides: slide- Chapter 1.ppt;
"@
$ws.Range("C3").Value = $c3Text
$ws.Rows.Item(3).RowHeight = 409.6

# --- C4: Week 2 topic (new content) ---
$c4Text = @"
Topic: name- Week 2, lectures- 2 Lectures, duration- 01:40;
Video: link- https://drive.google.com/file/d/15p5EP1_UCBHxyD9URakXw9xEW3hBKCLu/preview, name- TPL Week#2 Part 1, duration- 00:51;
Video: link- https://drive.google.com/file/d/1hHFe7EcX_PNGHntkNqR0a4Q9tTmLY1d-/preview, name- TPL Week#2 Part 2, duration- 00:49;
Important: 1 question will come from Chapter 2;
Slides: slide- Chapter 2.ppt;
Slides: slide- Chapter 3.ppt;
Topics: Readability, reliability, writeability, Chapter 3, Lexemes, Backus Naur Form (BNF), Parse Trees;
"@
$ws.Range("C4").Value = $c4Text
$ws.Rows.Item(4).RowHeight = 272

# --- C5: Week 3 topic incl. assignment (new content) ---
$c5Text = @"
Topic: name- Week 3, lectures- 1 Lecture, duration- 01:13;
Video: link- https://drive.google.com/file/d/1Pu_2p-pfLQAGEPfub2GY7qYzCmXxi8LG/preview, name- TPL Week#3, duration- 01:13;
Topics: EBNF, Attribute Grammer, EBNF Parse Trees;
Important: Question can come that in words explain BNF,;
Assignment: name- Assignment # 1 | Compare For Loops for Java with C++, img- /TPL/Assignment1.png;
AssignmentSolution: name- Assignment 1 Solution | Compare For Loops, link- Assignment No 1 TPL Waqas 3151.docx;
"@
$ws.Range("C5").Font.Bold = $true
$ws.Range("C5").WrapText = $true
$ws.Range("C5").Value = $c5Text
$ws.Rows.Item(5).RowHeight = 272

# --- C6: Week 4 topic (new content) ---
$c6Text = @"
Topic: name- Week 4, lectures- 1 Lecture, duration- 00:43;
Video: link- https://drive.google.com/file/d/1JZjwQWO7_aTq04qMJLL5pf6HXNrv_qyN/preview, name- TPL Week#4, duration- 00:43;
Slides: slide- Chapter 4.ppt;
"@
$ws.Range("C6").Font.Bold = $true
$ws.Range("C6").WrapText = $true
$ws.Range("C6").Value = $c6Text

# --- C7: Week 5 topic (new content) ---
$c7Text = @"
Topic: name- Week 5, lectures- 1 Lecture, duration- 01:20;
Video: link- https://drive.google.com/file/d/11P2oVuOQdNCB1lXC0azUnH74t3y6MwIx/preview, name- TPL Week#5, duration- 01:20;
Slides: slide- Chapter 5.ppt
"@
$ws.Range("C7").Font.Bold = $true
$ws.Range("C7").WrapText = $true
$ws.Range("C7").Value = $c7Text

# --- View state: active cell moved to D4, top-left scrolled to row 4 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("D4").Select()
